$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the "Main Menu" block (current rows 33-36) down one row, making
#    room for a new blank separator row at row 33 (the start of the new
#    "Game Over Screen" section). A true row insert correctly carries both
#    values and per-cell styles along with it.
# ---------------------------------------------------------------------------
[void]$ws.Rows("33:33").Insert(-4121, 1)

# ---------------------------------------------------------------------------
# 2. Build the new dark "section separator" fill off-sheet, in scratch
#    space, so we only pay the style-table cost of creating it once, then
#    stamp it onto the real target cells with a single format-only paste.
# ---------------------------------------------------------------------------
[void]$ws.Range("A29").Copy()
[void]$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1").Interior.ThemeColor = 1

[void]$ws.Range("B12").Copy()
[void]$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("Z2").Interior.ThemeColor = 1

[void]$ws.Range("Z1").Copy()
[void]$ws.Range("A33").PasteSpecial(-4122)

[void]$ws.Range("Z2").Copy()
[void]$ws.Range("B33:F33").PasteSpecial(-4122)

[void]$ws.Range("A33:F33").ClearContents()

# Reset the scratch cells back to an untouched, default-styled, empty state.
[void]$ws.Range("Z100").Copy()
[void]$ws.Range("Z1:Z2").PasteSpecial(-4122)
[void]$ws.Range("Z1:Z2").ClearContents()

# ---------------------------------------------------------------------------
# 3. Append the new "Game Over Screen" rows (38-39), matching the existing
#    bordered/filled look of the other data rows.
# ---------------------------------------------------------------------------
[void]$ws.Range("B12").Copy()
[void]$ws.Range("A38:F39").PasteSpecial(-4122)

$ws.Range("A38").Value = "Try Again button"
$ws.Range("B38").Value = 384
$ws.Range("C38").Value = 300
$ws.Range("D38").Value = 288
$ws.Range("E38").Value = 58
$ws.Range("F38").Value = "Game Over Screen"

$ws.Range("A39").Value = "Quit Game Over"
$ws.Range("B39").Value = 460
$ws.Range("C39").Value = 376
$ws.Range("D39").Value = 144
$ws.Range("E39").Value = 58
$ws.Range("F39").Value = "Game Over Screen"

# ---------------------------------------------------------------------------
# 4. Fix up sheet-level bookkeeping: view scroll position + selection.
# ---------------------------------------------------------------------------
[void]$ws.Range("A33:F33").Select()

"done"
